# Update "想去人数" (want-to-go count) values in the F column across the
# relevant worksheets, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets(1))
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 846
$wsExpo.Range("F6").Value  = 611
$wsExpo.Range("F13").Value = 1280
$wsExpo.Range("F28").Value = 44
$wsExpo.Range("F29").Value = 575
$wsExpo.Range("F31").Value = 245
$wsExpo.Range("F36").Value = 75

# Sheet "演出" (Worksheets(2))
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F18").Value = 64
$wsShow.Range("F22").Value = 284
$wsShow.Range("F28").Value = 8
$wsShow.Range("F37").Value = 647
$wsShow.Range("F38").Value = 647

# Sheet "全部类型" (Worksheets(4))
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value  = 846
$wsAll.Range("F10").Value = 611
$wsAll.Range("F20").Value = 1280
$wsAll.Range("F38").Value = 44
$wsAll.Range("F41").Value = 575
$wsAll.Range("F43").Value = 245
$wsAll.Range("F49").Value = 75
$wsAll.Range("F51").Value = 647
